$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$st = $ws.Range("D2").Style
$ws.Range("D2").Value = "'64.009.77"
$ws.Range("D2").Style = $st
$st = $ws.Range("E2").Style
$ws.Range("E2").Value = "'  +5.74%  "
$ws.Range("E2").Style = $st
$st = $ws.Range("D3").Style
$ws.Range("D3").Value = "'3.154.39"
$ws.Range("D3").Style = $st
$st = $ws.Range("E3").Style
$ws.Range("E3").Value = "'  +4.07%  "
$ws.Range("E3").Style = $st
$st = $ws.Range("E4").Style
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("E4").Style = $st
$st = $ws.Range("D5").Style
$ws.Range("D5").Value = "'592.10"
$ws.Range("D5").Style = $st
$st = $ws.Range("E5").Style
$ws.Range("E5").Value = "'  +3.59%  "
$ws.Range("E5").Style = $st
$st = $ws.Range("D6").Style
$ws.Range("D6").Value = "'148.05"
$ws.Range("D6").Style = $st
$st = $ws.Range("E6").Style
$ws.Range("E6").Value = "'  +4.79%  "
$ws.Range("E6").Style = $st
$st = $ws.Range("E7").Style
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = $st
$st = $ws.Range("D8").Style
$ws.Range("D8").Value = "'3.144.92"
$ws.Range("D8").Style = $st
$st = $ws.Range("E8").Style
$ws.Range("E8").Value = "'  +3.83%  "
$ws.Range("E8").Style = $st
$st = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.537"
$ws.Range("D9").Style = $st
$st = $ws.Range("E9").Style
$ws.Range("E9").Value = "'  +2.63%  "
$ws.Range("E9").Style = $st
$st = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.164"
$ws.Range("D10").Style = $st
$st = $ws.Range("E10").Style
$ws.Range("E10").Value = "'  +19.25%  "
$ws.Range("E10").Style = $st
$st = $ws.Range("E11").Style
$ws.Range("E11").Value = "'  +6.12%  "
$ws.Range("E11").Style = $st
$st = $ws.Range("E12").Style
$ws.Range("E12").Value = "'  +3.04%  "
$ws.Range("E12").Style = $st
$st = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("D13").Style = $st
$st = $ws.Range("E13").Style
$ws.Range("E13").Value = "'  +8.61%  "
$ws.Range("E13").Style = $st
$st = $ws.Range("D14").Style
$ws.Range("D14").Value = "'36.09"
$ws.Range("D14").Style = $st
$st = $ws.Range("E14").Style
$ws.Range("E14").Value = "'  +4.30%  "
$ws.Range("E14").Style = $st
$st = $ws.Range("E15").Style
$ws.Range("E15").Value = "'  +0.50%  "
$ws.Range("E15").Style = $st
$st = $ws.Range("D16").Style
$ws.Range("D16").Value = "'3.677.32"
$ws.Range("D16").Style = $st
$st = $ws.Range("E16").Style
$ws.Range("E16").Value = "'  +4.14%  "
$ws.Range("E16").Style = $st
$st = $ws.Range("D17").Style
$ws.Range("D17").Value = "'7.21"
$ws.Range("D17").Style = $st
$st = $ws.Range("E17").Style
$ws.Range("E17").Value = "'  +0.64%  "
$ws.Range("E17").Style = $st
$st = $ws.Range("D18").Style
$ws.Range("D18").Value = "'63.943.47"
$ws.Range("D18").Style = $st
$st = $ws.Range("E18").Style
$ws.Range("E18").Value = "'  +5.64%  "
$ws.Range("E18").Style = $st
$st = $ws.Range("D19").Style
$ws.Range("D19").Value = "'3.146.87"
$ws.Range("D19").Style = $st
$st = $ws.Range("E19").Style
$ws.Range("E19").Value = "'  +3.67%  "
$ws.Range("E19").Style = $st
$st = $ws.Range("D20").Style
$ws.Range("D20").Value = "'472.10"
$ws.Range("D20").Style = $st
$st = $ws.Range("E20").Style
$ws.Range("E20").Value = "'  +6.36%  "
$ws.Range("E20").Style = $st
$st = $ws.Range("D21").Style
$ws.Range("D21").Value = "'14.27"
$ws.Range("D21").Style = $st
$st = $ws.Range("E21").Style
$ws.Range("E21").Value = "'  +3.42%  "
$ws.Range("E21").Style = $st
$st = $ws.Range("D22").Style
$ws.Range("D22").Value = "'0.736"
$ws.Range("D22").Style = $st
$st = $ws.Range("E22").Style
$ws.Range("E22").Value = "'  +0.93%  "
$ws.Range("E22").Style = $st
$st = $ws.Range("D23").Style
$ws.Range("D23").Value = "'7.58"
$ws.Range("D23").Style = $st
$st = $ws.Range("E23").Style
$ws.Range("E23").Value = "'  +5.36%  "
$ws.Range("E23").Style = $st
$st = $ws.Range("D24").Style
$ws.Range("D24").Value = "'13.43"
$ws.Range("D24").Style = $st
$st = $ws.Range("E24").Style
$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("E24").Style = $st
$st = $ws.Range("D25").Style
$ws.Range("D25").Value = "'82.65"
$ws.Range("D25").Style = $st
$st = $ws.Range("E25").Style
$ws.Range("E25").Value = "'  +1.60%  "
$ws.Range("E25").Style = $st
$st = $ws.Range("E26").Style
$ws.Range("E26").Value = "'  +0.06%  "
$ws.Range("E26").Style = $st
$st = $ws.Range("E27").Style
$ws.Range("E27").Value = "'  +9.56%  "
$ws.Range("E27").Style = $st
$st = $ws.Range("D28").Style
$ws.Range("D28").Value = "'2.72"
$ws.Range("D28").Style = $st
$st = $ws.Range("E28").Style
$ws.Range("E28").Value = "'  +4.99%  "
$ws.Range("E28").Style = $st
$st = $ws.Range("E29").Style
$ws.Range("E29").Value = "'  -1.89%  "
$ws.Range("E29").Style = $st
$st = $ws.Range("E30").Style
$ws.Range("E30").Value = "'  +0.32%  "
$ws.Range("E30").Style = $st
$st = $ws.Range("D31").Style
$ws.Range("D31").Value = "'6.86"
$ws.Range("D31").Style = $st
$st = $ws.Range("E31").Style
$ws.Range("E31").Value = "'  +7.96%  "
$ws.Range("E31").Style = $st
$st = $ws.Range("D32").Style
$ws.Range("D32").Value = "'27.19"
$ws.Range("D32").Style = $st
$st = $ws.Range("E32").Style
$ws.Range("E32").Value = "'  +3.30%  "
$ws.Range("E32").Style = $st
$st = $ws.Range("E33").Style
$ws.Range("E33").Value = "'  +3.58%  "
$ws.Range("E33").Style = $st
$st = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.0₃0879"
$ws.Range("D34").Style = $st
$st = $ws.Range("E34").Style
$ws.Range("E34").Value = "'  +9.50%  "
$ws.Range("E34").Style = $st
$st = $ws.Range("E35").Style
$ws.Range("E35").Value = "'  +13.95%  "
$ws.Range("E35").Style = $st
$st = $ws.Range("E36").Style
$ws.Range("E36").Value = "'  +3.22%  "
$ws.Range("E36").Style = $st
$st = $ws.Range("E37").Style
$ws.Range("E37").Value = "'  +15.42%  "
$ws.Range("E37").Style = $st
$st = $ws.Range("D38").Style
$ws.Range("D38").Value = "'6.18"
$ws.Range("D38").Style = $st
$st = $ws.Range("E38").Style
$ws.Range("E38").Value = "'  +2.56%  "
$ws.Range("E38").Style = $st
$st = $ws.Range("D39").Style
$ws.Range("D39").Value = "'50.99"
$ws.Range("D39").Style = $st
$st = $ws.Range("E39").Style
$ws.Range("E39").Value = "'  +2.74%  "
$ws.Range("E39").Style = $st
$st = $ws.Range("D40").Style
$ws.Range("D40").Value = "'450.85"
$ws.Range("D40").Style = $st
$st = $ws.Range("E40").Style
$ws.Range("E40").Value = "'  +11.00%  "
$ws.Range("E40").Style = $st
$st = $ws.Range("E41").Style
$ws.Range("E41").Value = "'  +0.07%  "
$ws.Range("E41").Style = $st
$st = $ws.Range("E42").Style
$ws.Range("E42").Value = "'  +5.34%  "
$ws.Range("E42").Style = $st
$st = $ws.Range("D43").Style
$ws.Range("D43").Value = "'2.931.84"
$ws.Range("D43").Style = $st
$st = $ws.Range("E43").Style
$ws.Range("E43").Value = "'  +5.79%  "
$ws.Range("E43").Style = $st
$st = $ws.Range("E44").Style
$ws.Range("E44").Value = "'  +9.99%  "
$ws.Range("E44").Style = $st
$st = $ws.Range("E45").Style
$ws.Range("E45").Value = "'  +5.59%  "
$ws.Range("E45").Style = $st
$st = $ws.Range("E46").Style
$ws.Range("E46").Value = "'  +6.72%  "
$ws.Range("E46").Style = $st
$st = $ws.Range("D47").Style
$ws.Range("D47").Value = "'124.80"
$ws.Range("D47").Style = $st
$st = $ws.Range("E47").Style
$ws.Range("E47").Value = "'  +1.63%  "
$ws.Range("E47").Style = $st
$st = $ws.Range("E49").Style
$ws.Range("E49").Value = "'  +1.85%  "
$ws.Range("E49").Style = $st
$st = $ws.Range("D50").Style
$ws.Range("D50").Value = "'34.61"
$ws.Range("D50").Style = $st
$st = $ws.Range("E50").Style
$ws.Range("E50").Value = "'  -5.11%  "
$ws.Range("E50").Style = $st
$st = $ws.Range("D51").Style
$ws.Range("D51").Value = "'25.02"
$ws.Range("D51").Style = $st
$st = $ws.Range("E51").Style
$ws.Range("E51").Value = "'  +5.67%  "
$ws.Range("E51").Style = $st
